$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 546 (a new weekly price block for Plátano),
# pushing the existing rows 546:656 down to 549:659.
$ws.Rows("546:548").Insert()

# Populate the 3 newly inserted rows with the new weekly observation
# (date 2022-05-13 / serial 44694), reusing the same template values
# (market/category/quality) as the rest of the sheet.

# Row 546 - Pintón
$ws.Cells.Item(546,1).Value = 8
$ws.Cells.Item(546,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(546,3).Value = "Coquimbo"
$ws.Cells.Item(546,4).Value = 44694
$ws.Cells.Item(546,5).Value = 4
$ws.Cells.Item(546,6).Value = "Fruta"
$ws.Cells.Item(546,7).Value = 100108
$ws.Cells.Item(546,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(546,9).Value = 100108006
$ws.Cells.Item(546,10).Value = "Plátano"
$ws.Cells.Item(546,11).Value = "Sin especificar"
$ws.Cells.Item(546,12).Value = "Pintón"
$ws.Cells.Item(546,13).Value = 80
$ws.Cells.Item(546,14).Value = 11000
$ws.Cells.Item(546,15).Value = 11000
$ws.Cells.Item(546,16).Value = 11000
$ws.Cells.Item(546,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(546,18).Value = "Ecuador"
$ws.Cells.Item(546,19).Value = 550
$ws.Cells.Item(546,20).Value = 20

# Row 547 - Primera Maduro
$ws.Cells.Item(547,1).Value = 8
$ws.Cells.Item(547,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(547,3).Value = "Coquimbo"
$ws.Cells.Item(547,4).Value = 44694
$ws.Cells.Item(547,5).Value = 4
$ws.Cells.Item(547,6).Value = "Fruta"
$ws.Cells.Item(547,7).Value = 100108
$ws.Cells.Item(547,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(547,9).Value = 100108006
$ws.Cells.Item(547,10).Value = "Plátano"
$ws.Cells.Item(547,11).Value = "Sin especificar"
$ws.Cells.Item(547,12).Value = "Primera Maduro"
$ws.Cells.Item(547,13).Value = 120
$ws.Cells.Item(547,14).Value = 13000
$ws.Cells.Item(547,15).Value = 13000
$ws.Cells.Item(547,16).Value = 13000
$ws.Cells.Item(547,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(547,18).Value = "Ecuador"
$ws.Cells.Item(547,19).Value = 650
$ws.Cells.Item(547,20).Value = 20

# Row 548 - Primera Pintón
$ws.Cells.Item(548,1).Value = 8
$ws.Cells.Item(548,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(548,3).Value = "Coquimbo"
$ws.Cells.Item(548,4).Value = 44694
$ws.Cells.Item(548,5).Value = 4
$ws.Cells.Item(548,6).Value = "Fruta"
$ws.Cells.Item(548,7).Value = 100108
$ws.Cells.Item(548,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(548,9).Value = 100108006
$ws.Cells.Item(548,10).Value = "Plátano"
$ws.Cells.Item(548,11).Value = "Sin especificar"
$ws.Cells.Item(548,12).Value = "Primera Pintón"
$ws.Cells.Item(548,13).Value = 120
$ws.Cells.Item(548,14).Value = 14000
$ws.Cells.Item(548,15).Value = 14000
$ws.Cells.Item(548,16).Value = 14000
$ws.Cells.Item(548,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(548,18).Value = "Ecuador"
$ws.Cells.Item(548,19).Value = 700
$ws.Cells.Item(548,20).Value = 20

Write-Host "done"
